# The edit replaces the (English) "Note from Chuck ..." acknowledgement
# paragraph that lives in the *speaker notes* of slide 1 with its Greek
# translation. That paragraph is not on the slide body itself -- it is
# the text of the notes-page body placeholder shape attached to slide 1,
# reached through Slide.NotesPage.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$np = $s.NotesPage

$greek = "Σημείωση από τον  Chuck. Εάν χρησιμοποιείτε αυτό το υλικό, μπορείτε να αφαιρέσετε το λογότυπο UM και να το αντικαταστήσετε με το δικό σας, αλλά διατηρήστε το λογότυπο CC-BY στην πρώτη σελίδα καθώς την/τις σελίδα/ες αναγνώρισης."

# Find the notes-page shape that currently holds the "Note from Chuck"
# acknowledgement text (normally the first shape / body placeholder),
# and replace its text with the Greek translation.
$targetShape = $null
for ($i = 1; $i -le $np.Shapes.Count; $i++) {
    $candidate = $np.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.TextRange.Text -like "*Note from Chuck*") {
        $targetShape = $candidate
        break
    }
}
if ($targetShape -eq $null) {
    $targetShape = $np.Shapes.Item(1)
}

$targetShape.TextFrame.TextRange.Text = $greek
